$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12, shifting the existing rows 12-14 down to 13-15.
$ws.Rows.Item(12).EntireRow.Insert()

# Populate the newly inserted row 12 with the new weekly data point.
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value = "La Araucanía"
$ws.Cells.Item(12, 4).Value = 44813
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 100112036
$ws.Cells.Item(12, 7).Value = "Caigua"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 20
$ws.Cells.Item(12, 11).Value = 20000
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 20000
$ws.Cells.Item(12, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 1333
$ws.Cells.Item(12, 17).Value = 15
$ws.Cells.Item(12, 18).Value = "Hortaliza"
